# status: finish turn PM
# Applies the "studies" workbook edit:
#  - Rename sheet "estudos" -> "Estudos"
#  - Add a "DESCANSO" column (F) to Tabela1, extending the table to A1:F3
#  - Update ASSUNTO/PRODUCAO text for row 2 (rewording) and add same for row 3
#  - Fill in the previously-empty row 3 (HORA F, ASSUNTO, PRODUCAO) + new DESCANSO values
#  - Keep row 3 visually consistent with row 2 (height, wrapped text styling)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Estudos"

# 2. Add the new "DESCANSO" table column (this extends Tabela1 from A1:E3 to A1:F3)
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# 3. Header cell F1 - match the look of the other header cells (copy format from E1)
$ws.Range("E1").Copy()
$null = $ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "DESCANSO"

# 4. Row 2: ASSUNTO becomes "HARD", PRODUCAO gets the corrected/expanded wording
$ws.Range("D2").Value = "HARD"
$ws.Range("E2").Value = "Aula de HARD (Controle de sessão) + Início de implemetação do controle de produtos o (projeto green collections)"

# 5. Row 2: new DESCANSO (break) duration, 00:10 stored as a time fraction of a day
$ws.Range("F2").NumberFormat = "h:mm;@"
$ws.Range("F2").Value = 0.0069444444444444441

# 6. Row 3: HORA I cell formatting should follow the same look as the rest of the
#    time columns (copy the format from B2, which already carries that style)
$ws.Range("B2").Copy()
$null = $ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = 0.58819444444444446

# 7. Row 3: HORA F (end time)
$ws.Range("C3").NumberFormat = "h:mm;@"
$ws.Range("C3").Value = 0.76736111111111116

# 8. Row 3: ASSUNTO / PRODUCAO text
$ws.Range("D3").Value = "HARD"
$ws.Range("E3").Value = "Implementação de controle de produtos (projeto green collections) + Ajuda alpha (Atividade extra)"

# 9. Row 3: DESCANSO duration, 00:15
$ws.Range("F3").NumberFormat = "h:mm;@"
$ws.Range("F3").Value = 0.010416666666666666

# 10. Row 3 now holds wrapped text like row 2, so give it the same row height
$ws.Rows.Item(3).RowHeight = 31.5

# 11. Leave the selection/view where the author left off (scrolled right one column,
#     caret resting on the first empty row beneath the table)
$ws.Application.ActiveWindow.ScrollColumn = 2
$null = $ws.Range("C4").Select()
